# Lithuania A Lyga base update (24-02-2024 21:58)
#
# Two pairs of adjacent match rows had their B:AC data (id + all odds/stat
# columns) swapped between rows, while column A (the running row index) and
# the row/worksheet structure stay untouched:
#   - row 125 <-> row 126
#   - row 164 <-> row 165

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $dataA = $rangeA.Value2
    $dataB = $rangeB.Value2

    $rangeA.Value2 = $dataB
    $rangeB.Value2 = $dataA
}

Swap-RowData 125 126
Swap-RowData 164 165
